$wb = $excel.ActiveWorkbook

# --- Update the conversion text on "Hoja1" (sheet1), cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.52 = 5528.81 pesos`n✅ 5528.81 pesos = 1.52 = 862.87 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $text

# --- Update the rate values on "tasas" (sheet2) ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 659.999
$ws2.Range("O10").Value = 3649.01
$ws2.Range("N12").Value = 3645.85
$ws2.Range("O12").Value = 569
